$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated quarter-on-quarter confidence interval bounds (Low = column C, High = column D)
# following a change in the method used to compute the CIs. A couple of Tracker
# (column B) values also shifted by floating point noise from the recompute.
$ws.Range("C2").Value = 0.5473019062485873
$ws.Range("D2").Value = 8.288602944482349
$ws.Range("C3").Value = 2.336538197474081
$ws.Range("D3").Value = 8.420925159145698
$ws.Range("C4").Value = -5.831459815107165
$ws.Range("D4").Value = 0.6234400863243339
$ws.Range("C5").Value = -5.636901255877069
$ws.Range("D5").Value = 0.4363987909712375
$ws.Range("C6").Value = 1.860057215531952
$ws.Range("D6").Value = 6.55452694385259
$ws.Range("C7").Value = -3.469059361761806
$ws.Range("D7").Value = 4.044228171379527
$ws.Range("C8").Value = -2.501748859922448
$ws.Range("D8").Value = 1.715977296859039
$ws.Range("C9").Value = -1.026487832326983
$ws.Range("D9").Value = 8.07446689491873
$ws.Range("C10").Value = -1.891620692598783
$ws.Range("D10").Value = 4.678778177908094
$ws.Range("C11").Value = 4.697622100979482
$ws.Range("D11").Value = 10.23117127235673
$ws.Range("C12").Value = -5.868388739426789
$ws.Range("D12").Value = 2.149573713386155
$ws.Range("C13").Value = -3.397858862689174
$ws.Range("D13").Value = 2.202899183221585
$ws.Range("C14").Value = -2.247565637251847
$ws.Range("D14").Value = 9.582552649439414
$ws.Range("C15").Value = -2.093849368727063
$ws.Range("D15").Value = 4.486396289097594
$ws.Range("C16").Value = -8.206427647871639
$ws.Range("D16").Value = -2.388814416703777
$ws.Range("C17").Value = -5.781331925614486
$ws.Range("D17").Value = -0.6737226720217704
$ws.Range("C18").Value = -4.875271439088824
$ws.Range("D18").Value = 5.673794054999504
$ws.Range("B19").Value = -3.233103716856356
$ws.Range("C19").Value = -6.79272002661514
$ws.Range("D19").Value = 0.354039452534427
$ws.Range("C20").Value = -7.029876087379394
$ws.Range("D20").Value = 16.61525665687831
$ws.Range("C21").Value = 5.172767133968859
$ws.Range("D21").Value = 12.28708311906699
$ws.Range("C22").Value = -3.629395256345935
$ws.Range("D22").Value = 5.193432380973029
$ws.Range("C23").Value = -5.715341457707702
$ws.Range("D23").Value = 8.646068047429578
$ws.Range("B24").Value = 2.889754734408911
$ws.Range("C24").Value = -1.626334856703671
$ws.Range("D24").Value = 7.042899899225441
$ws.Range("C25").Value = -9.28758839008329
$ws.Range("D25").Value = -2.567162040239357
$ws.Range("C26").Value = -1.995882370427571
$ws.Range("D26").Value = 2.172343226982432
$ws.Range("C27").Value = -2.457443688369687
$ws.Range("D27").Value = 3.511402829487098
$ws.Range("C28").Value = -9.540190186306319
$ws.Range("D28").Value = 1.851338375843281
$ws.Range("C29").Value = -5.703349166695681
$ws.Range("D29").Value = 3.434529998689184
$ws.Range("C30").Value = -7.883192160846653
$ws.Range("D30").Value = 3.391976355012405
$ws.Range("C31").Value = -0.132097697872835
$ws.Range("D31").Value = 5.791697719361966
$ws.Range("C32").Value = -3.655847344489593
$ws.Range("D32").Value = 0.9938410023416999
$ws.Range("C33").Value = 1.900217324304676
$ws.Range("D33").Value = 7.994627897889983
$ws.Range("C34").Value = -2.061375672485499
$ws.Range("D34").Value = 3.801018057852512
$ws.Range("C35").Value = -6.472377563371811
$ws.Range("D35").Value = -1.362223302897791
$ws.Range("C36").Value = -3.42875879737371
$ws.Range("D36").Value = 3.137309320630766
$ws.Range("C37").Value = -2.991585692138232
$ws.Range("D37").Value = 4.115610394396385
$ws.Range("C38").Value = -3.676000334441809
$ws.Range("D38").Value = 2.969678772021767
$ws.Range("C39").Value = -5.06665337150568
$ws.Range("D39").Value = 3.949433207136255
$ws.Range("C40").Value = -10.8951310919454
$ws.Range("D40").Value = -0.3375281937973429
$ws.Range("C41").Value = 0.6901253132752361
$ws.Range("D41").Value = 6.589542588563235
$ws.Range("C42").Value = -1.721483162111392
$ws.Range("D42").Value = 4.101776676959679
$ws.Range("C43").Value = -1.294177313495837
$ws.Range("D43").Value = 4.051273380556597
$ws.Range("C44").Value = -4.323823421201311
$ws.Range("D44").Value = 0.8249507386989263
$ws.Range("C45").Value = -4.7945702258127
$ws.Range("D45").Value = 4.664635848788401
$ws.Range("C46").Value = -2.683594819994117
$ws.Range("D46").Value = 3.359661780826251
$ws.Range("C47").Value = -0.2157515231801521
$ws.Range("D47").Value = 4.836737369774324
